$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column J (rand_digit) on Sheet1, per row.
$values = @{
    2  = 2
    3  = 8
    4  = 3
    5  = 1
    6  = 8
    7  = 1
    8  = 2
    9  = 1
    10 = 5
    11 = 1
    12 = 2
    13 = 6
    14 = 8
    15 = 7
    16 = 3
    17 = 2
    18 = 4
    19 = 2
    20 = 4
    21 = 2
    22 = 8
    24 = 7
    25 = 6
    26 = 1
    27 = 2
    28 = 1
    29 = 6
    30 = 6
    31 = 1
    32 = 1
    33 = 8
    34 = 4
    35 = 5
    36 = 8
    37 = 2
    38 = 4
    39 = 3
    41 = 7
}

foreach ($row in $values.Keys) {
    $ws.Range("J$row").Value = $values[$row]
}
